$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
    @{ Row = 10; A = 42613.759398148148; B = 22 },
    @{ Row = 11; A = 42613.888055555559; B = 22 },
    @{ Row = 12; A = 42614.886203703703; B = 52 },
    @{ Row = 13; A = 42615.885671296295; B = 1 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = 0
    $ws.Range("D$row").Value = 0
    $ws.Range("E$row").Value = 0
    $ws.Range("F$row").Value = 0
    $ws.Range("G$row").Value = 0
    $ws.Range("H$row").Value = 0
    $ws.Range("I$row").Value = 0
    $ws.Range("J$row").Value = 0
    $ws.Range("K$row").Value = 0
    $ws.Range("L$row").Value = 0
    $ws.Range("M$row").Value = 0
    $ws.Range("N$row").Value = "Random"
}
